$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$text = 'questions = [
    {
        "title": "Your plant manufactures fiber-reinforced plastic materials, using resin as a primary raw material. You receive an equal amount of resin from two suppliers and must determine the stock reorder point for this material. Your plant''s daily resin consumption and your suppliers'' delivery time are outlined in the tables below.  What is the stock reorder point when you will place orders to both suppliers?",
        "ques_type": 2,
        "options": [
            "38,000 liters",
            "43,000 liters",
            "48,000 liters",
            "66,000 liters"
        ],
        "score": "43,000 liters"
    },
    {
        "title": "Your plant quality team has to decide on a single sampling technique for all received shipments, based on the Quality Control (QC) pass ratio Q. The sampling technique suitable for supplies received from most of your vendors will be selected as the sole sampling technique for your plant. Last month\u2019s gearbox shipments received from four suppliers are summarized in the table below. Which sampling technique should you select?",
        "ques_type": 2,
        "options": [
            " Simple random sampling",
            " Systematic sampling",
            " Complete testing",
            "Quota testing"
        ],
        "score": "Systematic sampling"
    },
    {
        "title": "You are analyzing mold manufacturing time to increase the throughput of your plant. A mold passes through machine type A &gt B &gt C &gt D in that sequence. Your plant initially had one machine of each type. The plant manager ordered an additional unit each for machines A and B. Still, they complain that the plant throughput has not increased. You review the table below summarizing machine count and processing time. Which machine is the bottleneck?",
        "ques_type": 2,
        "options": [
            "A",
            "B",
            "C",
            "D"
        ],
        "score": "C"
    },
    {
        "title": "Today\u2019s dispatch is to four customers in the south, identified in the first table. You\u2019re considering implementing a milk-run delivery replenishment model for these customers from closer warehouses instead of from the plant. Your logistics manager informed you that this would be financially beneficial only if: 1)     The total distance of customers from the warehouse is &lt500km.2)     A milk run serves no customers between the plant and the warehouse. From which warehouse(s) will you deploy a milk run?",
        "ques_type": 2,
        "options": [
            "Warehouse P",
            "Warehouse Q",
            "Both warehouses",
            "Neither warehouse"
        ],
        "score": "Warehouse Q"
    }
]'
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
